$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value for each data row (rows 2-173).
# Update all of these from 2023-09-01 (45170) to 2023-09-05 (45174).
$ws.Range("C2:C173").Value = 45174
